# "Added Week 15 simulations"
# Appends a new week's worth of simulated game data to the per-play yardage
# series (YDS sheet) and the special-teams distance series (ST sheet), and
# updates the cumulative season totals on OFF / DEF / ST / TURNS / PEN to
# reflect the newly simulated week.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append newly simulated per-play yardage numbers to each
# of the four running series (Offense-Rush, Offense-Pass, Defense-Rush,
# Defense-Pass).
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value2 = $ydsWs.Range("B2").Value2 + " 9 13 -1 3 -3 5 3 32 7 -3 8 5 9 3 13 -1 3 13 3 0 3 6 3 47 4 3 0 4 2"
$ydsWs.Range("B3").Value2 = $ydsWs.Range("B3").Value2 + " 1 29 5 8 24 25 16 55 7 13 22 14 1 29 6 1 4"
$ydsWs.Range("C2").Value2 = $ydsWs.Range("C2").Value2 + " -1 4 0 1 6 5 -6 2 10 0 2 2 -1 5 2 9 0 -3 0 8 3 5 -4 7 4 3"
$ydsWs.Range("C3").Value2 = $ydsWs.Range("C3").Value2 + " 7 17 13 30 5 10 7 0 15 8 11 15 15 12 6 9 18 7 11 -2 19 10 7 11 9 11 12 4 17 6 18 7 2 6 8"

# ---------------------------------------------------------------------
# OFF sheet: cumulative season totals (Home row2 / Road row3)
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value2 = 345
$offWs.Range("E2").Value2 = 21
$offWs.Range("F2").Value2 = 89
$offWs.Range("G2").Value2 = 96
$offWs.Range("H2").Value2 = 18
$offWs.Range("I2").Value2 = 11
$offWs.Range("J2").Value2 = 48

$offWs.Range("B3").Value2 = 13
$offWs.Range("C3").Value2 = 322
$offWs.Range("E3").Value2 = 62
$offWs.Range("F3").Value2 = 172
$offWs.Range("G3").Value2 = 84
$offWs.Range("I3").Value2 = 90
$offWs.Range("J3").Value2 = 96
$offWs.Range("L3").Value2 = 497
$offWs.Range("M3").Value2 = 332
$offWs.Range("Q3").Value2 = 881

# ---------------------------------------------------------------------
# DEF sheet: cumulative season totals (Home row2 / Road row3)
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value2 = 370
$defWs.Range("D2").Value2 = 24
$defWs.Range("F2").Value2 = 118
$defWs.Range("G2").Value2 = 118
$defWs.Range("I2").Value2 = 14
$defWs.Range("J2").Value2 = 67
$defWs.Range("N2").Value2 = 37
$defWs.Range("O2").Value2 = 48

$defWs.Range("C3").Value2 = 407
$defWs.Range("D3").Value2 = 11
$defWs.Range("E3").Value2 = 67
$defWs.Range("F3").Value2 = 226
$defWs.Range("G3").Value2 = 78
$defWs.Range("H3").Value2 = 48
$defWs.Range("I3").Value2 = 120
$defWs.Range("J3").Value2 = 99
$defWs.Range("L3").Value2 = 633
$defWs.Range("M3").Value2 = 438
$defWs.Range("Q3").Value2 = 1074

# ---------------------------------------------------------------------
# ST sheet: cumulative totals plus the four distance series (D / RA / RM
# each have two series: one in column B and one in column D).
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value2 = 146
$stWs.Range("D2").Value2 = 137
$stWs.Range("F2").Value2 = 250
$stWs.Range("G2").Value2 = 224
$stWs.Range("H2").Value2 = 5
$stWs.Range("I2").Value2 = 2
$stWs.Range("J2").Value2 = 87
$stWs.Range("K2").Value2 = 82

$stWs.Range("B3").Value2 = 69

$stWs.Range("D3").Value2 = $stWs.Range("D3").Value2 + " 41 48 48 32"
$stWs.Range("B4").Value2 = $stWs.Range("B4").Value2 + " 66 61 50"
$stWs.Range("D4").Value2 = $stWs.Range("D4").Value2 + " 0 2 0 0"
$stWs.Range("B5").Value2 = $stWs.Range("B5").Value2 + " 13 24 11"
$stWs.Range("D5").Value2 = "0 0 0 0 0 9 22"
$stWs.Range("B6").Value2 = $stWs.Range("B6").Value2 + " 15 18 25"

# ---------------------------------------------------------------------
# TURNS sheet: Road FMBL total increases
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("D3").Value2 = 14

# ---------------------------------------------------------------------
# PEN sheet: penalty counts increase
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value2 = 30
$penWs.Range("B3").Value2 = 27
$penWs.Range("D3").Value2 = 9
$penWs.Range("D4").Value2 = 12

# ---------------------------------------------------------------------
# Active tab moves from ST back to YDS
# ---------------------------------------------------------------------
$ydsWs.Activate()
